$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Hgf"
$ws.Cells.Item(2,3).Value = "Sdc2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 5.231719
$ws.Cells.Item(2,8).Value = 15.695157
$ws.Cells.Item(2,9).Value = 0.362499186434781
$ws.Cells.Item(2,10).Value = 0.362499186434781
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 1.321445333333333
$ws.Cells.Item(2,14).Value = 3.964336
$ws.Cells.Item(2,15).Value = 0.01021782062667047
$ws.Cells.Item(2,16).Value = 0.01021782062667047
$ws.Cells.Item(2,17).Value = 6.913430657861333
$ws.Cells.Item(2,18).Value = 62.22087592075201
$ws.Cells.Item(2,19).Value = 0.003703951664304568
$ws.Cells.Item(2,20).Value = 0.003703951664304569

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Hgf"
$ws.Cells.Item(3,3).Value = "Sdc2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 5.231719
$ws.Cells.Item(3,8).Value = 15.695157
$ws.Cells.Item(3,9).Value = 0.362499186434781
$ws.Cells.Item(3,10).Value = 0.362499186434781
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 105.9632263333333
$ws.Cells.Item(3,14).Value = 317.889679
$ws.Cells.Item(3,15).Value = 0.819340166699254
$ws.Cells.Item(3,16).Value = 0.8193401666992541
$ws.Cells.Item(3,17).Value = 554.3698245094004
$ws.Cells.Item(3,18).Value = 4989.328420584603
$ws.Cells.Item(3,19).Value = 0.2970101438418174
$ws.Cells.Item(3,20).Value = 0.2970101438418175

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Hgf"
$ws.Cells.Item(4,3).Value = "Sdc2"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 5.231719
$ws.Cells.Item(4,8).Value = 15.695157
$ws.Cells.Item(4,9).Value = 0.362499186434781
$ws.Cells.Item(4,10).Value = 0.362499186434781
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 22.04284166666666
$ws.Cells.Item(4,14).Value = 66.128525
$ws.Cells.Item(4,15).Value = 0.1704420126740755
$ws.Cells.Item(4,16).Value = 0.1704420126740755
$ws.Cells.Item(4,17).Value = 115.3219535614917
$ws.Cells.Item(4,18).Value = 1037.897582053425
$ws.Cells.Item(4,19).Value = 0.06178509092865898
$ws.Cells.Item(4,20).Value = 0.06178509092865899

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Hgf"
$ws.Cells.Item(5,3).Value = "Sdc2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 8.316945333333335
$ws.Cells.Item(5,8).Value = 24.950836
$ws.Cells.Item(5,9).Value = 0.5762706133406404
$ws.Cells.Item(5,10).Value = 0.5762706133406403
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 1.321445333333333
$ws.Cells.Item(5,14).Value = 3.964336
$ws.Cells.Item(5,15).Value = 0.01021782062667047
$ws.Cells.Item(5,16).Value = 0.01021782062667047
$ws.Cells.Item(5,17).Value = 10.99038859832178
$ws.Cells.Item(5,18).Value = 98.91349738489602
$ws.Cells.Item(5,19).Value = 0.005888229759536037
$ws.Cells.Item(5,20).Value = 0.005888229759536037

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Hgf"
$ws.Cells.Item(6,3).Value = "Sdc2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 8.316945333333335
$ws.Cells.Item(6,8).Value = 24.950836
$ws.Cells.Item(6,9).Value = 0.5762706133406404
$ws.Cells.Item(6,10).Value = 0.5762706133406403
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 105.9632263333333
$ws.Cells.Item(6,14).Value = 317.889679
$ws.Cells.Item(6,15).Value = 0.819340166699254
$ws.Cells.Item(6,16).Value = 0.8193401666992541
$ws.Cells.Item(6,17).Value = 881.2903607579607
$ws.Cells.Item(6,18).Value = 7931.613246821645
$ws.Cells.Item(6,19).Value = 0.4721616603984017
$ws.Cells.Item(6,20).Value = 0.4721616603984017

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Hgf"
$ws.Cells.Item(7,3).Value = "Sdc2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 8.316945333333335
$ws.Cells.Item(7,8).Value = 24.950836
$ws.Cells.Item(7,9).Value = 0.5762706133406404
$ws.Cells.Item(7,10).Value = 0.5762706133406403
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 22.04284166666666
$ws.Cells.Item(7,14).Value = 66.128525
$ws.Cells.Item(7,15).Value = 0.1704420126740755
$ws.Cells.Item(7,16).Value = 0.1704420126740755
$ws.Cells.Item(7,17).Value = 183.3291091329889
$ws.Cells.Item(7,18).Value = 1649.9619821969
$ws.Cells.Item(7,19).Value = 0.09822072318270268
$ws.Cells.Item(7,20).Value = 0.09822072318270268

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Hgf"
$ws.Cells.Item(8,3).Value = "Sdc2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.8836963333333333
$ws.Cells.Item(8,8).Value = 2.651089
$ws.Cells.Item(8,9).Value = 0.06123020022457864
$ws.Cells.Item(8,10).Value = 0.06123020022457864
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 1.321445333333333
$ws.Cells.Item(8,14).Value = 3.964336
$ws.Cells.Item(8,15).Value = 0.01021782062667047
$ws.Cells.Item(8,16).Value = 0.01021782062667047
$ws.Cells.Item(8,17).Value = 1.167756395767111
$ws.Cells.Item(8,18).Value = 10.509807561904
$ws.Cells.Item(8,19).Value = 0.0006256392028298624
$ws.Cells.Item(8,20).Value = 0.0006256392028298625

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Hgf"
$ws.Cells.Item(9,3).Value = "Sdc2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.8836963333333333
$ws.Cells.Item(9,8).Value = 2.651089
$ws.Cells.Item(9,9).Value = 0.06123020022457864
$ws.Cells.Item(9,10).Value = 0.06123020022457864
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 105.9632263333333
$ws.Cells.Item(9,14).Value = 317.889679
$ws.Cells.Item(9,15).Value = 0.819340166699254
$ws.Cells.Item(9,16).Value = 0.8193401666992541
$ws.Cells.Item(9,17).Value = 93.63931457893678
$ws.Cells.Item(9,18).Value = 842.7538312104309
$ws.Cells.Item(9,19).Value = 0.05016836245903496
$ws.Cells.Item(9,20).Value = 0.05016836245903497

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Hgf"
$ws.Cells.Item(10,3).Value = "Sdc2"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.8836963333333333
$ws.Cells.Item(10,8).Value = 2.651089
$ws.Cells.Item(10,9).Value = 0.06123020022457864
$ws.Cells.Item(10,10).Value = 0.06123020022457864
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 22.04284166666666
$ws.Cells.Item(10,14).Value = 66.128525
$ws.Cells.Item(10,15).Value = 0.1704420126740755
$ws.Cells.Item(10,16).Value = 0.1704420126740755
$ws.Cells.Item(10,17).Value = 19.47917835708055
$ws.Cells.Item(10,18).Value = 175.312605213725
$ws.Cells.Item(10,19).Value = 0.01043619856271381
$ws.Cells.Item(10,20).Value = 0.01043619856271381

Write-Output "Edit applied: rows 2-10 updated for FAPs/sCs/ECs x Hgf/Sdc2 pairs"
